$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A81").NumberFormat = "@"
$ws.Range("A81").Value = "2026/01/30"
$ws.Range("B81").Value = "逃离鸭科夫"
$ws.Range("C81").Value = 1167
$ws.Range("A80:C80").Copy()
$ws.Range("A81:C81").PasteSpecial(-4122)
